# "dependencies fix and test" - correct the attendance counts (P.Days / T.Days)
# on Sheet1.  Columns: A=SN., B=Name, C=Attn, D=Date, E=P.Days, F=T.Days.
#   Row 2 (Bill Gates): T.Days 3 -> 2
#   Row 3 (GMTK Guy):   P.Days 3 -> 1, T.Days 3 -> 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 2
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
